# Auto-generated edit script applying value updates described in the diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("G4").Value = 70
$ws.Range("F5").Value = 482
$ws.Range("F7").Value = 1150
$ws.Range("F9").Value = 178
$ws.Range("F11").Value = 767
$ws.Range("F12").Value = 412
$ws.Range("F15").Value = 203
$ws.Range("F18").Value = 6247
$ws.Range("F22").Value = 7220
$ws.Range("F25").Value = 3294
$ws.Range("F26").Value = 420
$ws.Range("F27").Value = 797
$ws.Range("F31").Value = 158
$ws.Range("F32").Value = 1289
$ws.Range("F33").Value = 121
$ws.Range("F36").Value = 1005
$ws.Range("F37").Value = 1312
$ws.Range("F38").Value = 2078

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 61

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 61
$ws.Range("G7").Value = 70
$ws.Range("F8").Value = 482
$ws.Range("F10").Value = 1150
$ws.Range("F12").Value = 178
$ws.Range("F14").Value = 767
$ws.Range("F15").Value = 412
$ws.Range("F19").Value = 203
$ws.Range("F22").Value = 6247
$ws.Range("F23").Value = 6247
$ws.Range("F27").Value = 7220
$ws.Range("F30").Value = 3294
$ws.Range("F31").Value = 420
$ws.Range("F32").Value = 797
$ws.Range("F37").Value = 158
$ws.Range("F38").Value = 1289
$ws.Range("F39").Value = 121
$ws.Range("F42").Value = 1005
$ws.Range("F43").Value = 1312
$ws.Range("F45").Value = 2078
